$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in password for performance_glitch_user row (B6): s1ecret_sauce -> secret_sauce
$ws.Range("B6").Value = "secret_sauce"

# Update the active selection to E10
$ws.Range("E10").Select()
